# Restructure the header on the KSRO, KSRM and KSPO sheets:
#   - the merged A1:A2 / B1:K1 header pair becomes a single A1:K1 merge
#   - "BSSID STATUS" moves up into the new merged header (A1)
#   - "Floor" moves down into A2, above the per-column location labels

$wb = $excel.ActiveWorkbook

$sheetNames = @("KSRO", "KSRM", "KSPO")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Break apart the two existing merges so the cells can be edited freely.
    $ws.Range("A1:A2").UnMerge()
    $ws.Range("B1:K1").UnMerge()

    # Swap the header text: "BSSID STATUS" now spans the header row,
    # "Floor" becomes the label above the location row.
    $ws.Range("A1").Value = "BSSID STATUS"
    $ws.Range("B1").Value = $null
    $ws.Range("A2").Value = "Floor"

    # Recombine into a single header merge across the whole row.
    $ws.Range("A1:K1").Merge()

    # Merging recalculates borders per-cell; restore the original uniform
    # box-border style (as used by every other cell) by re-pasting the
    # formatting from an untouched cell in the same style family.
    $ws.Range("C2").Copy()
    $ws.Range("A1:K1").PasteSpecial(-4122) | Out-Null
    $ws.Range("A2").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}
